$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @(-0.1362899550126131, 2.08867529575408,  20.09818370085882,  4.48309978707354,  4.581737469302817, 23)
    3  = @(1.479266720037939,  2.107761023088795,  21.2282041533365,   4.607407530633306, 4.466166091718527, 22)
    4  = @(0.1953051818482703, 1.330132550180737,  4.813685238236514,  2.194011221082635, 2.239267331157384, 21)
    5  = @(0.7065468732113371, 0.8830325039473388, 1.138482200557635,  1.066996813752335, 0.8203168300068852, 20)
    6  = @(0.6543568374413384, 0.8911237843702812, 1.363273558160643,  1.167593061884423, 0.9934989532409025, 19)
    7  = @(0.2055713410900566, 0.6219053496127894, 0.5362322189693115, 0.732278785005623, 0.723208031597916,  18)
    8  = @(0.2066150860819048, 0.5311179486106835, 0.416429203296134,  0.6453132598173805, 0.6301572205356532, 17)
    9  = @(0.2632867251892846, 0.2910766374824335, 0.1134320491269004, 0.3367967475004775, 0.220280195234621,  11)
    10 = @(0.2847860619040521, 0.2847860619040521, 0.1099954283298057, 0.3316555869117926, 0.1835966098656279, 7)
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    $ws.Cells.Item($row, 2).Value = $rowVals[0]
    $ws.Cells.Item($row, 3).Value = $rowVals[1]
    $ws.Cells.Item($row, 4).Value = $rowVals[2]
    $ws.Cells.Item($row, 5).Value = $rowVals[3]
    $ws.Cells.Item($row, 6).Value = $rowVals[4]
    $ws.Cells.Item($row, 7).Value = $rowVals[5]
}
